$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B20: "cor-5000" -> "Cor-500" (new distinct shared string)
$ws.Range("B20").Value = "Cor-500"

# Rename the existing shared string used by B6 / B15 from "cor-5000" -> "cor-500"
$ws.Range("B6").Value = "cor-500"
$ws.Range("B15").Value = "cor-500"

# New execution-time values for RF / cor-500..L1000-tm (rows 6-10)
$ws.Range("C6").Value = 3411.4386
$ws.Range("C7").Value = 2693.9418
$ws.Range("C8").Value = 2693.9418
$ws.Range("C9").Value = 3553.4004
$ws.Range("C10").Value = 3411.4386

# New execution-time values for EN / text-mining..L1000-tm (rows 14-19)
$ws.Range("C14").Value = 121.6083
$ws.Range("C15").Value = 76.89216
$ws.Range("C16").Value = 61.52664
$ws.Range("C17").Value = 78.564
$ws.Range("C18").Value = 74.66988
$ws.Range("C19").Value = 93.49674

# Match the updated view/selection state (A14->A1 scroll, A17->C5 active cell)
$ws.Range("C5").Select() | Out-Null
